$d = $word.ActiveDocument

# --- Fix 1: "About Project: ... Jquery)" -- merge the split runs around
#     "Jquery" back into a single run and drop the spell-check proofErr
#     markers that bracketed it.
$d.Content.Find.Execute(
    "About Project: An E-commerce website making (using Core PHP/ajax/Jquery)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "About Project: An E-commerce website making (using Core PHP/ajax/Jquery)",
    2) | Out-Null

# --- Fix 2: "Project Name: Asbab Furniture" -- merge "Asbab" and
#     " Furniture" (and drop their proofErr wrapper) while keeping the
#     red colour only on "Asbab Furniture", not on "Project Name: ".
$d.Content.Find.Execute(
    "Project Name: Asbab Furniture",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Project Name: Asbab Furniture",
    2) | Out-Null

$colorRng = $d.Content
$colorRng.Find.Execute("Asbab Furniture", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($colorRng.Find.Found) {
    $colorRng.Font.Color = 255
}

# --- Fix 3: table cell "6" (week 6, column 1) gains extra blank-ish
#     paragraphs: trailing spaces on the "6" line, a line of spaces, and
#     a "done" marker line, matching the pattern used by the other week
#     rows in this table.
$table = $d.Tables.Item(1)
$targetCell = $null
for ($row = 1; $row -le $table.Rows.Count; $row++) {
    $cellText = $table.Cell($row, 1).Range.Text
    if ($cellText.TrimEnd([char]13, [char]7) -eq "6") {
        $targetCell = $table.Cell($row, 1)
        break
    }
}

if ($targetCell -ne $null) {
    $cellPara = $targetCell.Range.Paragraphs.Item(1)
    $insertPoint = $cellPara.Range.Duplicate
    $insertPoint.MoveEnd(1, -1) | Out-Null
    $insertPoint.Collapse(0)

    $cr = [char]13
    $newText = "       " + $cr + "             " + $cr + "                                                                        " + "d" + "one"
    $insertPoint.InsertAfter($newText) | Out-Null
}

Write-Host "edit applied"
